# hs-logger_priorities.xlsx — mark "Change explicit list access to using .get()"
# as complete: move it out of the in-progress block (row 5) down into the
# not-started block, landing just above the "Generate graphs automatically"
# row, and flip its Todo flag (column B) from 1 to 0. The Priority column
# (E) is a formula (=B*C/D) so it recalculates to 0 on its own.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 5 holds the task being completed; capture its text + Importance/Difficulty
# before we shuffle rows around.
$taskName   = $ws.Cells.Item(5, 1).Text
$importance = $ws.Cells.Item(5, 3).Value2
$difficulty = $ws.Cells.Item(5, 4).Value2

# Remove row 5 (rows 6:81 shift up to 5:80), then open a fresh row at 35
# (rows 35:80 shift back down to 36:81) and drop the task there instead.
$ws.Rows(5).Delete()
$ws.Rows(35).Insert()

$ws.Cells.Item(35, 1).Value = $taskName
$ws.Cells.Item(35, 2).Value = 0
$ws.Cells.Item(35, 3).Value = $importance
$ws.Cells.Item(35, 4).Value = $difficulty
$ws.Cells.Item(35, 5).Formula = "=B35*C35/D35"

# The conditional-formatting rules on column B covered B3:B80; the table
# still ends at row 81, so extend those two color-scale rules to match.
$oldRange = $ws.Range("B3:B80")
$newRange = $ws.Range("B3:B81")
$fcs = $oldRange.FormatConditions
for ($i = 1; $i -le $fcs.Count; $i++) {
    $fc = $fcs.Item($i)
    if ($fc.AppliesTo.Address() -eq $oldRange.Address()) {
        $fc.ModifyAppliesToRange($newRange)
    }
}

# Leave the active selection where the author left it.
$ws.Range("E15").Select()
